$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) column cells being updated so that
# values such as "0.0540" or "4.50" keep their trailing zeros as text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.656.07"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "2.643.75"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").Value = "146.48"
$ws.Range("E6").Value = "  +3.81%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").Value = "6.92"
$ws.Range("E9").Value = "  +7.21%  "
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "3.112.79"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "59.557.38"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "21.44"
$ws.Range("E15").Value = "  +4.26%  "
$ws.Range("D16").Value = "2.629.27"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D18").Value = "4.50"
$ws.Range("E18").Value = "  +3.01%  "
$ws.Range("D19").Value = "340.05"
$ws.Range("D20").Value = "10.33"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "66.45"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").Value = "0.418"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").Value = "0.0₃0750"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "1.65"
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("D31").Value = "5.87"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").Value = "18.87"
$ws.Range("D33").Value = "151.07"
$ws.Range("D34").Value = "4.01"
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("D36").Value = "0.839"
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("D37").Value = "0.840"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("D40").Value = "285.89"
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "0.604"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "0.0540"
$ws.Range("E43").Value = "  +3.14%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "10.74"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "19.25"
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("D46").Value = "0.0946"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").Value = "1.964.88"
$ws.Range("D49").Value = "18.50"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").Value = "4.57"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").Value = "111.50"
$ws.Range("E51").Value = "  +0.15%  "
